# Goblin_Profits market-data refresh (scheduled runner)
# Updates cached currentAveragePrice / leve profit columns (H:N) across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 597.7407
$ws.Range("I28").Value = 551.2381
$ws.Range("J28").Value = 760.5
$ws.Range("K28").Value = 551.2381
$ws.Range("L28").Value = 760.5
$ws.Range("M28").Value = -66.23810000000003
$ws.Range("N28").Value = -1730.5
$ws.Range("H100").Value = 6248.8
$ws.Range("I100").Value = 3496.3333
$ws.Range("K100").Value = 3496.3333
$ws.Range("M100").Value = -2955.3333
$ws.Range("H111").Value = 600.5
$ws.Range("I111").Value = 527.5714
$ws.Range("J111").Value = 702.6
$ws.Range("K111").Value = 1582.7142
$ws.Range("L111").Value = 2107.8
$ws.Range("M111").Value = 1484.2858
$ws.Range("N111").Value = -8241.8
$ws.Range("H132").Value = 4002162.0
$ws.Range("I132").Value = 2245.1904
$ws.Range("K132").Value = 6735.5712
$ws.Range("M132").Value = -4205.5712

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2075.9565
$ws.Range("I45").Value = 1689.875
$ws.Range("J45").Value = 2958.4285
$ws.Range("K45").Value = 1689.875
$ws.Range("L45").Value = 2958.4285
$ws.Range("M45").Value = -1312.875
$ws.Range("N45").Value = -3712.4285
$ws.Range("H74").Value = 2411.3403
$ws.Range("I74").Value = 2082.7585
$ws.Range("J74").Value = 2940.7222
$ws.Range("K74").Value = 2082.7585
$ws.Range("L74").Value = 2940.7222
$ws.Range("M74").Value = -1208.7585
$ws.Range("N74").Value = -4688.7222
$ws.Range("H77").Value = 2411.3403
$ws.Range("I77").Value = 2082.7585
$ws.Range("J77").Value = 2940.7222
$ws.Range("K77").Value = 10413.7925
$ws.Range("L77").Value = 14703.611
$ws.Range("M77").Value = -6045.7925
$ws.Range("N77").Value = -23439.611

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 17498.625
$ws.Range("I26").Value = 17498.625
$ws.Range("J26").Value = 0.0
$ws.Range("K26").Value = 17498.625
$ws.Range("L26").Value = 0.0
$ws.Range("M26").Value = -17206.625
$ws.Range("N26").ClearContents()
$ws.Range("H100").Value = 20799.111
$ws.Range("J100").Value = 20799.111
$ws.Range("L100").Value = 20799.111
$ws.Range("N100").Value = -22963.111
$ws.Range("H134").Value = 863626.44
$ws.Range("I134").Value = 1189.409
$ws.Range("K134").Value = 3568.227
$ws.Range("M134").Value = -1033.227

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3563.3333
$ws.Range("I16").Value = 845.0
$ws.Range("K16").Value = 845.0
$ws.Range("M16").Value = -558.0
$ws.Range("H113").Value = 3563.3333
$ws.Range("I113").Value = 845.0
$ws.Range("K113").Value = 845.0
$ws.Range("M113").Value = 1325.0
$ws.Range("H122").Value = 30305134.0
$ws.Range("J122").Value = 41667850.0
$ws.Range("L122").Value = 125003550.0
$ws.Range("N122").Value = -125008450.0
$ws.Range("H132").Value = 5754.9165
$ws.Range("I132").Value = 1510.2858
$ws.Range("J132").Value = 11697.4
$ws.Range("K132").Value = 4530.857400000001
$ws.Range("L132").Value = 35092.2
$ws.Range("M132").Value = -2000.857400000001
$ws.Range("N132").Value = -40152.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 104.21429
$ws.Range("J2").Value = 162.375
$ws.Range("L2").Value = 974.25
$ws.Range("N2").Value = -1200.25
$ws.Range("H4").Value = 31823676.0
$ws.Range("I4").Value = 2317454.8
$ws.Range("K4").Value = 6952364.399999999
$ws.Range("M4").Value = -6952252.399999999
$ws.Range("H46").Value = 6900413.0
$ws.Range("I46").Value = 11616362.0
$ws.Range("J46").Value = 2858171.2
$ws.Range("K46").Value = 34849086.0
$ws.Range("L46").Value = 8574513.600000001
$ws.Range("M46").Value = -34848995.0
$ws.Range("N46").Value = -8574695.600000001
$ws.Range("H131").Value = 3336761.2
$ws.Range("J131").Value = 4448701.5
$ws.Range("L131").Value = 13346104.5
$ws.Range("N131").Value = -13356184.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 62504536.0
$ws.Range("I80").Value = 125002270.0
$ws.Range("J80").Value = 6800.25
$ws.Range("K80").Value = 125002270.0
$ws.Range("L80").Value = 6800.25
$ws.Range("M80").Value = -125001272.0
$ws.Range("N80").Value = -8796.25
$ws.Range("H83").Value = 62504536.0
$ws.Range("I83").Value = 125002270.0
$ws.Range("J83").Value = 6800.25
$ws.Range("K83").Value = 625011350.0
$ws.Range("L83").Value = 34001.25
$ws.Range("M83").Value = -625006358.0
$ws.Range("N83").Value = -43985.25
$ws.Range("H113").Value = 10909.0
$ws.Range("I113").Value = 9999.0
$ws.Range("K113").Value = 9999.0
$ws.Range("M113").Value = -7829.0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4145.636
$ws.Range("I7").Value = 4113.353
$ws.Range("J7").Value = 4255.4
$ws.Range("K7").Value = 4113.353
$ws.Range("L7").Value = 4255.4
$ws.Range("M7").Value = -4001.353
$ws.Range("N7").Value = -4479.4
$ws.Range("H44").Value = 9500.0
$ws.Range("J44").Value = 9500.0
$ws.Range("L44").Value = 9500.0
$ws.Range("N44").Value = -10412.0
$ws.Range("H82").Value = 52634056.0
$ws.Range("I82").Value = 728.1539
$ws.Range("J82").Value = 166672930.0
$ws.Range("K82").Value = 728.1539
$ws.Range("L82").Value = 166672930.0
$ws.Range("M82").Value = -367.1539
$ws.Range("N82").Value = -166673652.0
$ws.Range("H85").Value = 52634056.0
$ws.Range("I85").Value = 728.1539
$ws.Range("J85").Value = 166672930.0
$ws.Range("K85").Value = 728.1539
$ws.Range("L85").Value = 166672930.0
$ws.Range("M85").Value = 519.8461
$ws.Range("N85").Value = -166675426.0
$ws.Range("H105").Value = 50000.0
$ws.Range("J105").Value = 50000.0
$ws.Range("L105").Value = 50000.0
$ws.Range("M105").Value = -56988.0
$ws.Range("H122").Value = 4637.0
$ws.Range("I122").Value = 4517.727
$ws.Range("K122").Value = 13553.181
$ws.Range("M122").Value = -11103.181
$ws.Range("H126").Value = 4145.636
$ws.Range("I126").Value = 4113.353
$ws.Range("J126").Value = 4255.4
$ws.Range("K126").Value = 12340.059
$ws.Range("L126").Value = 12766.2
$ws.Range("M126").Value = -9870.059000000001
$ws.Range("N126").Value = -17706.2
$ws.Range("H136").Value = 17504.182
$ws.Range("I136").Value = 4974.3335
$ws.Range("K136").Value = 14923.0005
$ws.Range("M136").Value = -12373.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 606979.44
$ws.Range("I122").Value = 1266872.0
$ws.Range("J122").Value = 7077.091
$ws.Range("K122").Value = 3800616.0
$ws.Range("L122").Value = 21231.273
$ws.Range("M122").Value = -3798166.0
$ws.Range("N122").Value = -26131.273
